$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GmailInfo")

$ws.Range("A1").Value = "Critical security alert"

$body = "Sign-in attempt was blocked`nmailfortesting086@gmail.com`nSomeone just used your password to try to sign in to your account. Google blocked them, but you should check what happened.`nCheck activity`nYou received this email to let you know about important changes to your Google Account and services.`n© 2020 Google LLC, 1600 Amphitheatre Parkway, Mountain View, CA 94043, USA"

$ws.Range("B1").Value = $body
